$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the component_3_flowrate column (old column E); this shifts F:J left to E:I
$ws.Range("E:E").Delete()

# Update header labels and cell values/descriptions for the binary (2-component) system
$ws.Range("C2").Value = "component_A_flowrate"
$ws.Range("D2").Value = "component_B_flowrate"
$ws.Range("D3").Value = 60
$ws.Range("I3").Value = "Raw binary feed to distillation column"
$ws.Range("C4").Value = 95
$ws.Range("D4").Value = 2
$ws.Range("I4").Value = "Top product from column (A-rich)"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 58
$ws.Range("I5").Value = "Bottom product from column (B-rich)"
$ws.Range("C6").Value = 47.5
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 40
$ws.Range("B7").Value = "Net Distillate"
$ws.Range("C7").Value = 47.5
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 40
$ws.Range("I7").Value = "Final distillate product (A-rich)"
$ws.Range("B8").Value = "Reboiler Vapor"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 58
$ws.Range("H8").Value = "DISTILLATION_COL"
$ws.Range("I8").Value = "Vapor back to column"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 58
$ws.Range("E9").Value = 110
$ws.Range("I9").Value = "Final bottom product (B-rich)"
